$d = $word.ActiveDocument

# Paragraph 1 stays "Archivo 1 del repositorio repo" (the _GoBack bookmark that
# currently sits at its end will be relocated below).
$p1 = $d.Paragraphs(1)

# Add a new, empty paragraph right after paragraph 1.
$p1.Range.InsertParagraphAfter()

# Add a second new paragraph after that empty one; this will hold the new text.
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs(3)
$r3 = $p3.Range

# Insert "Modificación", remember the position right after it (the bookmark
# goes there), then insert the remaining " del archivo 1" text.
$insertStart = $r3.Start
$r3.InsertAfter("Modificación")
$bookmarkPos = $insertStart + "Modificación".Length
$r3.InsertAfter(" del archivo 1")

# Re-create the _GoBack bookmark between the two runs of paragraph 3; since
# _GoBack is a single-instance bookmark, this also removes it from paragraph 1.
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))
